# Update Reichelt Part List
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quantity of the "1x20pol female header" (row 9) increased from 1 to 3
$ws.Range("C9").Value = 3

# Price/piece (F9) now populated with the per-unit price; match the
# currency number format already used by the rest of column F/G
$ws.Range("F9").NumberFormat = $ws.Range("G9").NumberFormat
$ws.Range("F9").Value = 0.26

# Price all (G9) recalculated for the new quantity (3 * 0.26)
$ws.Range("G9").Value = 0.78

# Scroll / selection state: window was scrolled down so row 7 is the
# top-most visible row, and the active cell/selection moved to E10
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("E10").Select()

# Workbook window was resized (maximized) in the source edit
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 12435
$excel.ActiveWindow.WindowState = -4137
